# Applies the "fixed the bishop function and added results from tests" edit:
#  - Moves the existing "Can Game Piece Move Properly" secondary table
#    (AI10:AP17) down/left to (P27:W34)
#  - Moves the lone "Black"/"White" marker cells AL8 -> S25 and AL19 -> S36
#  - Adds a brand new test-results block (AF6:AK9) for the new
#    "Can Move To can't hit more than one enemy" unit test
#  - Updates the sheet view's active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- move the AI10:AP17 block to P27:W34 (shift: -19 cols, +17 rows) ---
$ws.Range("AI10:AP17").Copy($ws.Range("P27"))
$ws.Range("AI10:AP17").ClearContents()

# --- move the lone marker cells AL8 -> S25 and AL19 -> S36 ---
$ws.Range("AL8").Copy($ws.Range("S25"))
$ws.Range("AL8").ClearContents()

$ws.Range("AL19").Copy($ws.Range("S36"))
$ws.Range("AL19").ClearContents()

# --- new test block: "Can Move To can't hit more than one enemy" ---
$ws.Range("AF6").Value = "TESTS"
$ws.Range("AF7").Value = "Can Move To can't hit more than one enemy "

$ws.Range("AF8").Value = "Piece"
$ws.Range("AG8").Value = "Enemy 1"
$ws.Range("AH8").Value = "enemy 2"
$ws.Range("AI8").Value = "Desired Location"
$ws.Range("AJ8").Value = "Expected"
$ws.Range("AK8").Value = "Actual"

$ws.Range("AF9").Value = "White Bishop [D5]"
$ws.Range("AG9").Value = "E6"
$ws.Range("AH9").Value = "F7"
$ws.Range("AI9").Value = "F7"
$ws.Range("AJ9").Value = $false
$ws.Range("AK9").Value = "FAIL"
$ws.Range("AK9").Interior.Color = 255

# --- column widths for the two newly used columns ---
$ws.Columns("AF").ColumnWidth = 18.45
$ws.Columns("AI").ColumnWidth = 15.02

# --- update the sheet view selection to match the edited area ---
$ws.Range("U1").Select()
$ws.Range("AM16").Select()
